$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are stored as plain text in the source data (e.g. "592.92",
# "0.820", "65.039.84"). Force text format before assigning so Excel does not
# auto-convert numeric-looking values (which would also strip trailing zeros).
$priceCells = @("D2","D3","D5","D6","D7","D9","D11","D13","D14","D17","D18","D20","D22","D24","D25","D32","D33","D34","D39","D40","D42","D45","D47","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.039.84"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.516.82"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "592.92"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "133.95"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "3.517.95"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "7.14"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "4.117.94"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").Value = "27.61"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "3.521.51"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "65.036.15"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "14.32"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("D22").Value = "392.63"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "3.662.60"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "74.56"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("E28").Value = "  +8.25%  "
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "8.30"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "3.526.66"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "24.08"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").Value = "6.95"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "167.92"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "0.820"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "25.63"
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "4.44"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "2.418.94"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("E51").Value = "  +4.14%  "

# Restore the default (unstyled) cell style now that the text values are locked in,
# so the written cells match the original workbook formatting.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
